$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-9 from 2023-09-06 (45175)
# to 2023-09-14 (45183), keeping the existing date formatting/style.
for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 3).Value = 45183
}
